# The species-observation rows for ids 130965861 and 130965935 were
# re-ordered upstream: row 11 now holds what used to be row 12's data
# (the "Garnlav" / Alectoria sarmentosa record) and vice versa (the
# "Fläcknycklar" / Dactylorhiza maculata record). Swap the values of the
# columns that actually differ between the two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "AC", "AX")

foreach ($col in $cols) {
    $cell11 = $ws.Range($col + "11")
    $cell12 = $ws.Range($col + "12")

    $val11 = $cell11.Value2
    $val12 = $cell12.Value2

    $cell11.Value = $val12
    $cell12.Value = $val11
}
